$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-26 04:16:28"
$wsOverview.Range("G5").Value = "2016-08-26 04:16:28"

$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-26 04:16:24"
$wsZhCn.Range("H5").Value = "2016-08-26 04:16:24"
$wsZhCn.Range("K4").Value = "2016-08-26 04:16:40"
$wsZhCn.Range("K5").Value = "2016-08-26 04:16:40"

$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-26 04:16:28"
$wsDeDe.Range("H5").Value = "2016-08-26 04:16:28"
$wsDeDe.Range("K4").Value = "2016-08-26 04:16:47"
$wsDeDe.Range("K5").Value = "2016-08-26 04:16:47"
